# Swap the data content of rows 13 and 14 (records for "Spillkråka" /
# Dryocopus martius and "Dropptaggsvamp" / Hydnellum ferrugineum), keeping
# the shared/common columns (P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AY)
# untouched since they are identical between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($rng, [string]$val) {
    # Force text storage (these columns hold numeric-looking strings like
    # "1"/"5" as plain text, not as numbers) without leaving a lasting
    # style/number-format footprint on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# ---- Row 13 -> gets the former row-14 record ----
$ws.Range("A13").Value = 130881366
$ws.Range("B13").Value = 57881
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 100049
$ws.Range("F13").Value = "Spillkråka"
$ws.Range("G13").Value = "Dryocopus martius"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
Set-TextCell $ws.Range("I13") "1"
$ws.Range("J13").ClearContents()
Set-TextCell $ws.Range("K13") "adult"
Set-TextCell $ws.Range("M13") "äldre spår"
$ws.Range("Q13").Value = 656781
$ws.Range("R13").Value = 6559672
$ws.Range("Z13").Value = "08:48"
$ws.Range("AB13").Value = "08:48"
$ws.Range("AX13").Value = "Stuart Fell"

# ---- Row 14 -> gets the former row-13 record ----
$ws.Range("A14").Value = 130882201
$ws.Range("B14").Value = 93095
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 4364
$ws.Range("F14").Value = "Dropptaggsvamp"
$ws.Range("G14").Value = "Hydnellum ferrugineum"
$ws.Range("H14").Value = "(Fr.:Fr.) P. Karst."
Set-TextCell $ws.Range("I14") "5"
Set-TextCell $ws.Range("J14") "fruktkroppar"
$ws.Range("K14").ClearContents()
$ws.Range("M14").ClearContents()
$ws.Range("Q14").Value = 656955
$ws.Range("R14").Value = 6559350
$ws.Range("Z14").Value = "08:18"
$ws.Range("AB14").Value = "08:18"
$ws.Range("AX14").Value = "Stuart Fell, Liam Martin"
